# Update "Laundry Website Data Dictionary.xlsx"
# Commit message: "Update Laundry Website Data Dictionary.xlsx - Added all database variables"
#
# The data-dictionary table on Sheet1 lists database fields in columns:
#   B = P/F, C = Field Name, D = Caption, E = Data Type, F = Field Size, G = Notes
#
# Row 12 (userid) had its Data Type corrected from "varchar" to "varchar(20)".
# Rows 14 and 15 (previously blank) were filled in with two new fields:
#   week_day -> char, size 1, "Must be 1-7"
#   slot     -> slot, size 1, "Must be 1-8"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the Data Type for the existing "userid" row.
$ws.Range("E12").Value = "varchar(20)"

# Fill in the new "week_day" field row.
$ws.Range("C14").Value = "week_day"
$ws.Range("E14").Value = "char"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "Must be 1-7"

# Fill in the new "slot" field row.
$ws.Range("C15").Value = "slot"
$ws.Range("E15").Value = "slot"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "Must be 1-8"

# Update the active cell selection left by the editor.
$ws.Range("C17").Select()
